$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block for "Apio" (Femacal de La Calera) is a weekly rolling
# window of observations. This update drops the oldest observation
# (row 638, dated 2021-02-04) and appends a newly observed week at the
# end of the block (which lands at row 663 after the deletion), keeping
# every other row's position (and the rest of the sheet below it) intact.

# 1) Remove the oldest row - shifts rows 639..694 up to 638..693.
$ws.Rows.Item(638).Delete()

# 2) Re-open a blank row at 663 - shifts rows 663..693 back down to
#    664..694, restoring all rows that come after our block to their
#    original row numbers.
$ws.Rows.Item(663).Insert()

# 3) Populate the freshly inserted row with the new observation.
$ws.Cells.Item(663, 1).Value = 3
$ws.Cells.Item(663, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(663, 3).Value = "Coquimbo"
$ws.Cells.Item(663, 4).Value = 45147
$ws.Cells.Item(663, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(663, 5).Value = 5
$ws.Cells.Item(663, 6).Value = 100112017
$ws.Cells.Item(663, 7).Value = "Apio"
$ws.Cells.Item(663, 8).Value = "Americana (o)"
$ws.Cells.Item(663, 9).Value = "Primera"
$ws.Cells.Item(663, 10).Value = 110
$ws.Cells.Item(663, 11).Value = 8000
$ws.Cells.Item(663, 12).Value = 8000
$ws.Cells.Item(663, 13).Value = 8000
$ws.Cells.Item(663, 14).Value = "`$/docena de matas"
$ws.Cells.Item(663, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(663, 16).Value = 1333
$ws.Cells.Item(663, 17).Value = 6
$ws.Cells.Item(663, 18).Value = "Hortaliza"
